# The table rows for the following course groups had their "nazev"/"zkratka"
# shared-string entries swapped/rotated between rows (the idno in column C
# stays put per group):
#   Rows 7 & 8   (idno 4746): Mikroprocesory a senzory v praxi II/I <-> I/II
#   Rows 11 & 12 (idno 5232): zkratka ANE <-> KAEL (nazev stays "Analogová elektronika")
#   Rows 14,15,16 (idno 8514): nazev+zkratka rotate across the three rows
#   Rows 17 & 18 (idno 8753): zkratka OONV <-> KOONV (nazev stays "Objektově orientované návrhové vzory")
#
# Some zkratka values ("0164", "0171") look numeric, so a plain .Value=
# assignment would make Excel coerce them to numbers (losing the leading
# zero and the shared-string text type). To preserve them as literal text
# exactly like the original cells (t="s"), we stage the value in a scratch
# cell pre-formatted as Text, then copy only the VALUE (not formatting) onto
# the target cell via PasteSpecial, so the destination cell's style is left
# completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1000")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $scratch.Clear() | Out-Null
}

# Rows 7 & 8: swap nazev (A) and zkratka (B); idno (C) unchanged
$ws.Range("A7").Value = "Mikroprocesory a senzory v praxi I"
Set-TextValue $ws.Range("B7") "0164"
$ws.Range("A8").Value = "Mikroprocesory a senzory v praxi II"
Set-TextValue $ws.Range("B8") "0171"

# Rows 11 & 12: swap zkratka (B) only
$ws.Range("B11").Value = "KAEL"
$ws.Range("B12").Value = "ANE"

# Rows 14, 15, 16: rotate nazev (A) and zkratka (B)
$ws.Range("A14").Value = "Algoritmizace a programování I"
$ws.Range("B14").Value = "APR1"
$ws.Range("A15").Value = "Python and R for Data Science"
$ws.Range("B15").Value = "EPYR"
$ws.Range("A16").Value = "Algoritmizace a programování II"
$ws.Range("B16").Value = "APR2"

# Rows 17 & 18: swap zkratka (B) only
$ws.Range("B17").Value = "KOONV"
$ws.Range("B18").Value = "OONV"
